$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")
$summary = $wb.Worksheets.Item("Summary")

# Excel's Range.Value coerces numeric-looking strings ("10", "0146...") into
# real numbers (losing leading zeros / precision for long digit strings).
# The source workbook stores these as text, so force text via NumberFormat,
# then restore the cell's original style to avoid leaving a "@" format behind.
function Set-TextValue($range, [string]$value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# F21 changes from "1" to "10" (must remain text, not become a number)
Set-TextValue $ws.Range("F21") "10"

# New rows 22-31
$ws.Range("C22").Value = "506_紫罗兰香槟色_violet champagne_undefined_1bunch"
Set-TextValue $ws.Range("F22") "5"

Set-TextValue $ws.Range("A23") "4"
$ws.Range("C23").Value = "8_冰淇淋洋桔梗_Icecream Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
Set-TextValue $ws.Range("F23") "10"

$ws.Range("C24").Value = "14_波浪浅紫洋桔梗_Wavy Light Purple Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
Set-TextValue $ws.Range("F24") "10"

$ws.Range("C25").Value = "12_肉粉洋桔梗_Peach Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
Set-TextValue $ws.Range("F25") "15"

$ws.Range("C26").Value = "542_吊米 红_hanging amaranthus`nred_undefined_1bunch"
Set-TextValue $ws.Range("F26") "5"
$ws.Rows.Item(26).AutoFit()

$ws.Range("C27").Value = "448_吊米 绿_hanging amaranthus`ngreen_undefined_1bunch"
Set-TextValue $ws.Range("F27") "5"
$ws.Rows.Item(27).AutoFit()

$ws.Range("C28").Value = "630_吸色康乃馨天蓝_tinted tiffany blue_undefined_20stems"
Set-TextValue $ws.Range("F28") "10"

$ws.Range("C29").Value = "277_草莓杏仁饼_undefined_Rosa rugosa Thunb._10stems"
Set-TextValue $ws.Range("F29") "7"

$ws.Range("C30").Value = "480_蝴蝶洋牡丹红_butterfly  Ranunculus_undefined_1bunch"
Set-TextValue $ws.Range("F30") "5"

$ws.Range("C31").Value = "842_蝴蝶洋牡丹鲑鱼粉_undefined_undefined_1bunch"

# Summary sheet G2: long numeric-looking id, must stay exact text (no float rounding)
Set-TextValue $summary.Range("G2") "0146137101398786310151510155101051010155510750"
